# Convert the systematic/statistical uncertainty values in columns F and G
# (rows 2-3) from text labels such as "3.1*10^(-6)" into real numeric values,
# matching the author's re-upload of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.0000031
$ws.Range("G2").Value = 0.000003
$ws.Range("F3").Value = 0.0000064
$ws.Range("G3").Value = 0.0000031

# Match the saved selection/active cell from the re-uploaded file.
$ws.Range("H14").Select()
